$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.610.72"
$ws.Range("E2").Value = "'  +4.46%  "
$ws.Range("D3").Value = "'1.793.02"
$ws.Range("E3").Value = "'  +0.81%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.20%  "
$ws.Range("D5").Value = "'313.64"
$ws.Range("E5").Value = "'  -0.05%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "'  +0.38%  "
$ws.Range("D7").Value = "'0.5366"
$ws.Range("E7").Value = "'  +0.47%  "
$ws.Range("D8").Value = "'0.3809"
$ws.Range("E8").Value = "'  +1.07%  "
$ws.Range("D9").Value = "'0.07519"
$ws.Range("E9").Value = "'  +1.52%  "
$ws.Range("D10").Value = "'42.49"
$ws.Range("E10").Value = "'  -0.78%  "
$ws.Range("D11").Value = "'1.117"
$ws.Range("E11").Value = "'  +2.15%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "'  +0.15%  "
$ws.Range("D13").Value = "'21.08"
$ws.Range("E13").Value = "'  +1.84%  "
$ws.Range("D14").Value = "'6.176"
$ws.Range("E14").Value = "'  +1.24%  "
$ws.Range("D15").Value = "'7.388"
$ws.Range("E15").Value = "'  +5.65%  "
$ws.Range("D16").Value = "'1.792.98"
$ws.Range("E16").Value = "'  +0.59%  "
$ws.Range("D17").Value = "'90.33"
$ws.Range("E17").Value = "'  +0.76%  "
$ws.Range("D18").Value = "'0.00001066"
$ws.Range("E18").Value = "'  +0.97%  "
$ws.Range("D19").Value = "'0.06435"
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "'  +0.38%  "
$ws.Range("D21").Value = "'17.29"
$ws.Range("E21").Value = "'  +2.97%  "
$ws.Range("D22").Value = "'5.917"
$ws.Range("E22").Value = "'  +0.13%  "
$ws.Range("D23").Value = "'28.609.13"
$ws.Range("E23").Value = "'  +4.31%  "
$ws.Range("D24").Value = "'11.22"
$ws.Range("E24").Value = "'  +0.16%  "
$ws.Range("D25").Value = "'2.107"
$ws.Range("E25").Value = "'  +0.72%  "
$ws.Range("D26").Value = "'160.79"
$ws.Range("E26").Value = "'  +3.47%  "
$ws.Range("D27").Value = "'20.51"
$ws.Range("E27").Value = "'  +1.46%  "
$ws.Range("B28").Value = "'LidoDAOToken"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.364"
$ws.Range("E28").Value = "'  -0.35%  "
$ws.Range("B29").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "'2.002.02"
$ws.Range("E29").Value = "'  +0.69%  "
$ws.Range("D30").Value = "'123.40"
$ws.Range("E30").Value = "'  +1.74%  "
$ws.Range("D31").Value = "'1.121"
$ws.Range("E31").Value = "'  +3.74%  "
$ws.Range("D32").Value = "'0.1016"
$ws.Range("E32").Value = "'  -1.83%  "
$ws.Range("D33").Value = "'5.685"
$ws.Range("E33").Value = "'  +1.82%  "
$ws.Range("D34").Value = "'3.658"
$ws.Range("E34").Value = "'  +0.97%  "
$ws.Range("E35").Value = "'  +11.37%  "
$ws.Range("D36").Value = "'0.06577"
$ws.Range("E36").Value = "'  +10.20%  "
$ws.Range("D37").Value = "'0.02323"
$ws.Range("E37").Value = "'  +2.75%  "
$ws.Range("D38").Value = "'8.701"
$ws.Range("E38").Value = "'  +5.64%  "
$ws.Range("D39").Value = "'5.086"
$ws.Range("E39").Value = "'  +3.36%  "
$ws.Range("D40").Value = "'11.48"
$ws.Range("E40").Value = "'  +1.82%  "
$ws.Range("D41").Value = "'0.6322"
$ws.Range("E41").Value = "'  +3.19%  "
$ws.Range("E42").Value = "'  +6.36%  "
$ws.Range("D43").Value = "'1.003"
$ws.Range("E43").Value = "'  +0.36%  "
$ws.Range("D44").Value = "'1.379"
$ws.Range("E44").Value = "'  -3.42%  "
$ws.Range("D45").Value = "'13.52"
$ws.Range("E45").Value = "'  +2.29%  "
$ws.Range("D46").Value = "'0.5929"
$ws.Range("E46").Value = "'  +2.51%  "
$ws.Range("D47").Value = "'3.671"
$ws.Range("E47").Value = "'  +1.46%  "
$ws.Range("D48").Value = "'124.97"
$ws.Range("E48").Value = "'  +2.94%  "
$ws.Range("D49").Value = "'1.980"
$ws.Range("E49").Value = "'  +4.51%  "
$ws.Range("D50").Value = "'1.159"
$ws.Range("E50").Value = "'  +3.57%  "
$ws.Range("D51").Value = "'0.06920"
$ws.Range("E51").Value = "'  +2.89%  "
